$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (not auto-converted numbers) for numeric-looking price strings,
# matching the source workbook where these are plain text cells.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '72.361.42'
$ws.Range('E2').Value = '  +0.81%  '
$ws.Range('D3').Value = '4.050.59'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '533.52'
$ws.Range('E5').Value = '  +1.87%  '
$ws.Range('D6').Value = '152.15'
$ws.Range('E6').Value = '  +2.12%  '
$ws.Range('D7').Value = '0.700'
$ws.Range('E7').Value = '  +11.96%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '0.753'
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('D10').Value = '0.172'
$ws.Range('E10').Value = '  -3.17%  '
$ws.Range('D11').Value = '0.0000329'
$ws.Range('E11').Value = '  -3.53%  '
$ws.Range('D12').Value = '48.10'
$ws.Range('E12').Value = '  +3.64%  '
$ws.Range('D13').Value = '10.72'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('D14').Value = '4.663.74'
$ws.Range('D15').Value = '4.015.91'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('D16').Value = '14.24'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '20.71'
$ws.Range('E17').Value = '  -3.45%  '
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('D19').Value = '1.20'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').Value = '72.126.70'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('D21').Value = '429.71'
$ws.Range('E21').Value = '  -2.77%  '
$ws.Range('D22').Value = '98.67'
$ws.Range('E22').Value = '  +4.33%  '
$ws.Range('D23').Value = '3.50'
$ws.Range('E23').Value = '  -2.97%  '
$ws.Range('D24').Value = '4.23'
$ws.Range('E24').Value = '  +4.33%  '
$ws.Range('D25').Value = '14.38'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').Value = '11.25'
$ws.Range('E26').Value = '  -8.29%  '
$ws.Range('D27').Value = '10.81'
$ws.Range('E27').Value = '  -3.81%  '
$ws.Range('D28').Value = '5.85'
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('D29').Value = '36.94'
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('D30').Value = '3.59'
$ws.Range('E30').Value = '  +22.56%  '
$ws.Range('D31').Value = '13.44'
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').Value = '7.19'
$ws.Range('E32').Value = '  +3.21%  '
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('D34').Value = '677.21'
$ws.Range('E34').Value = '  -3.13%  '
$ws.Range('D35').Value = '44.74'
$ws.Range('E35').Value = '  +9.52%  '
$ws.Range('D36').Value = '66.06'
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('D37').Value = '0.452'
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0842'
$ws.Range('E38').Value = '  -7.69%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.153'
$ws.Range('E39').Value = '  -2.81%  '
$ws.Range('D40').Value = '3.40'
$ws.Range('E40').Value = '  -4.87%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = '0.0488'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').Value = '3.20'
$ws.Range('E44').Value = '  +2.37%  '
$ws.Range('D45').Value = '0.151'
$ws.Range('E45').Value = '  +3.20%  '
$ws.Range('D46').Value = '3.47'
$ws.Range('E46').Value = '  -1.45%  '
$ws.Range('D47').Value = '9.74'
$ws.Range('E47').Value = '  +5.97%  '
$ws.Range('D48').Value = '2.64'
$ws.Range('E48').Value = '  -6.12%  '
$ws.Range('D49').Value = '3.01'
$ws.Range('E49').Value = '  -6.08%  '
$ws.Range('D50').Value = '0.000271'
$ws.Range('E50').Value = '  -2.85%  '
$ws.Range('D51').Value = '145.31'
$ws.Range('E51').Value = '  +1.08%  '
